$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert one new column at B (old B..F "FI/ID/EX/MEM/WB" slide to
#    C..G). Column A (row labels) is untouched.
# ------------------------------------------------------------------
$ws.Range("B1").EntireColumn.Insert()

# ------------------------------------------------------------------
# 2. New headers: A1 = "inst", B1 = "width"
# ------------------------------------------------------------------
$ws.Range("A1").Value = "inst"
$ws.Range("B1").Value = "width"

# ------------------------------------------------------------------
# 3. New column B: bit-width values for each pipeline register row
# ------------------------------------------------------------------
$ws.Range("B2").Value = 32
$ws.Range("B3").Value = 32
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 4
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 32
$ws.Range("B11").Value = 32
$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 32
$ws.Range("B15").Value = 32

# ------------------------------------------------------------------
# 4. Mark which pipeline stage columns (now C:G = FI/ID/EX/MEM/WB)
#    carry each signal, by writing 1 into the relevant cells.
# ------------------------------------------------------------------
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1

# ------------------------------------------------------------------
# 5. Formatting.
#    - Header row (row 1) and the label column (A) keep the original
#      "output" look, light-grey fill, general/text format.
#    - The data grid (B2:G15) gets the new custom numeric format and a
#      uniform light-grey fill...
#    - ...except the "checked" cells, which are additionally
#      highlighted yellow.
# ------------------------------------------------------------------
$ws.Range("A1:G1").Interior.Color = 15921906
$ws.Range("A2:A15").Interior.Color = 15921906
$ws.Range("A1:G1").NumberFormat = "General"
$ws.Range("A2:A15").NumberFormat = "General"

$ws.Range("B2:G15").Interior.Color = 15921906
$ws.Range("B2:G15").NumberFormat = "0.00_);[Red]\(0.00\)"

$ws.Range("A1:G15").HorizontalAlignment = -4131
$ws.Range("A1:G15").VerticalAlignment = -4160

$checkedCells = @("C2","D2","C3","D3","D4","E4","F4","G4","D5","D6","D7","E7","D8","E8","F8","D9","E9","F9","D10","E10","D11","E11","F11","D12","E12","F12","G12","D13","E13","F13","G13","E14","F14","G14","F15","G15")
foreach ($addr in $checkedCells) {
    $ws.Range($addr).Interior.Color = 65535
}
